$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = 1127265073.931482
    $ws.Cells.Item($r, 5).Value = 1117.950775121333
    $ws.Cells.Item($r, 6).Value = -93432772.57103774
}
